# Update the "想去人数" (want-to-go count) column (F) values across the
# four worksheets to match the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 215
$ws1.Range("F4").Value  = 603
$ws1.Range("F6").Value  = 396
$ws1.Range("F7").Value  = 535
$ws1.Range("F12").Value = 567
$ws1.Range("F14").Value = 1718
$ws1.Range("F16").Value = 1714
$ws1.Range("F17").Value = 231
$ws1.Range("F18").Value = 480
$ws1.Range("F20").Value = 119

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value  = 216

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value  = 5252
$ws3.Range("F3").Value  = 295

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5252
$ws4.Range("F4").Value  = 295
$ws4.Range("F6").Value  = 215
$ws4.Range("F7").Value  = 216
$ws4.Range("F12").Value = 603
$ws4.Range("F16").Value = 396
$ws4.Range("F17").Value = 535
$ws4.Range("F25").Value = 567
$ws4.Range("F28").Value = 1718
$ws4.Range("F30").Value = 1716
$ws4.Range("F32").Value = 231
$ws4.Range("F33").Value = 480
$ws4.Range("F35").Value = 119

$wb.Save()
